$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prevent Excel from auto-converting date-like text (Y, AA columns) into date serials
$ws.Range("Y5:Y11").NumberFormat = "@"
$ws.Range("AA5:AA11").NumberFormat = "@"

# Row 5
$ws.Range("A5").Value = 74691146
$ws.Range("B5").Value = 89392
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 1202
$ws.Range("F5").Value = "Ullticka"
$ws.Range("G5").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H5").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P5").Value = "skogarna kring Holmtjärnsmossen, Holmtjärn och Sörgårdarna, Vstm"
$ws.Range("Q5").Value = 489647.9100812326
$ws.Range("R5").Value = 6631282.907482301
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = "Örebro"
$ws.Range("U5").Value = "Hällefors"
$ws.Range("V5").Value = "Västmanland"
$ws.Range("W5").Value = "Hjulsjö"
$ws.Range("Y5").Value = "2018-10-03"
$ws.Range("Z5").Value = "00:00"
$ws.Range("AA5").Value = "2018-10-03"
$ws.Range("AB5").Value = "00:00"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AI5").Value = "barrnaturskog"
$ws.Range("AO5").Value = "låga av gran"
$ws.Range("AW5").Value = "Sebastian Kirppu"
$ws.Range("AX5").Value = "Sebastian Kirppu"
$ws.Range("AY5").ClearContents()

# Row 6
$ws.Range("A6").Value = 82883467
$ws.Range("B6").Value = 95525
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 221941
$ws.Range("F6").Value = "Plattlummer"
$ws.Range("G6").Value = "Lycopodium complanatum"
$ws.Range("H6").Value = "L."
$ws.Range("P6").Value = "Holmtjärnen, Vstm"
$ws.Range("Q6").Value = 489447.1889090774
$ws.Range("R6").Value = 6631386.069839812
$ws.Range("S6").Value = 25
$ws.Range("T6").Value = "Örebro"
$ws.Range("U6").Value = "Hällefors"
$ws.Range("V6").Value = "Västmanland"
$ws.Range("W6").Value = "Hjulsjö"
$ws.Range("Y6").Value = "2019-10-17"
$ws.Range("Z6").Value = "00:00"
$ws.Range("AA6").Value = "2019-11-15"
$ws.Range("AB6").Value = "00:00"
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AI6").Value = "Äldre tallskog, lavtyp"
$ws.Range("AW6").Value = "Henrik Josefsson"
$ws.Range("AX6").Value = "Tommy Pettersson"
$ws.Range("AY6").Value = "Länsstyrelsen i Örebro län, inventering"

# Row 7
$ws.Range("A7").Value = 82883420
$ws.Range("B7").Value = 90645
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 4361
$ws.Range("F7").Value = "Orange taggsvamp"
$ws.Range("G7").Value = "Hydnellum aurantiacum"
$ws.Range("H7").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("P7").Value = "Holmtjärnen, Vstm"
$ws.Range("Q7").Value = 489475.2312818346
$ws.Range("R7").Value = 6631504.175627257
$ws.Range("S7").Value = 25
$ws.Range("T7").Value = "Örebro"
$ws.Range("U7").Value = "Hällefors"
$ws.Range("V7").Value = "Västmanland"
$ws.Range("W7").Value = "Hjulsjö"
$ws.Range("Y7").Value = "2019-10-17"
$ws.Range("Z7").Value = "00:00"
$ws.Range("AA7").Value = "2019-11-15"
$ws.Range("AB7").Value = "00:00"
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AI7").Value = "Äldre barrskog"
$ws.Range("AW7").Value = "Henrik Josefsson"
$ws.Range("AX7").Value = "Tommy Pettersson"
$ws.Range("AY7").Value = "Länsstyrelsen i Örebro län, inventering"

# Row 8
$ws.Range("A8").Value = 82883414
$ws.Range("B8").Value = 90697
$ws.Range("C8").Value = "Ovaliderad"
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 5449
$ws.Range("F8").Value = "Svart taggsvamp"
$ws.Range("G8").Value = "Phellodon niger"
$ws.Range("H8").Value = "(Fr.:Fr.) P.Karst."
$ws.Range("P8").Value = "Holmtjärnen, Vstm"
$ws.Range("Q8").Value = 489442.8430217677
$ws.Range("R8").Value = 6631274.934289954
$ws.Range("S8").Value = 25
$ws.Range("T8").Value = "Örebro"
$ws.Range("U8").Value = "Hällefors"
$ws.Range("V8").Value = "Västmanland"
$ws.Range("W8").Value = "Hjulsjö"
$ws.Range("Y8").Value = "2019-10-17"
$ws.Range("Z8").Value = "00:00"
$ws.Range("AA8").Value = "2019-11-15"
$ws.Range("AB8").Value = "00:00"
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AI8").Value = "Äldre barrskog, lav-ristyp"
$ws.Range("AO8").ClearContents()
$ws.Range("AW8").Value = "Henrik Josefsson"
$ws.Range("AX8").Value = "Tommy Pettersson"
$ws.Range("AY8").Value = "Länsstyrelsen i Örebro län, inventering"

# Row 9
$ws.Range("A9").Value = 82883445
$ws.Range("B9").Value = 94764
$ws.Range("C9").Value = "Ovaliderad"
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 2326
$ws.Range("F9").Value = "Vanlig rörsvepemossa"
$ws.Range("G9").Value = "Liochlaena lanceolata"
$ws.Range("H9").Value = "Nees"
$ws.Range("P9").Value = "Holmtjärnen, Vstm"
$ws.Range("Q9").Value = 489568.0120379663
$ws.Range("R9").Value = 6631356.056138818
$ws.Range("S9").Value = 25
$ws.Range("T9").Value = "Örebro"
$ws.Range("U9").Value = "Hällefors"
$ws.Range("V9").Value = "Västmanland"
$ws.Range("W9").Value = "Hjulsjö"
$ws.Range("Y9").Value = "2019-10-17"
$ws.Range("Z9").Value = "00:00"
$ws.Range("AA9").Value = "2019-11-15"
$ws.Range("AB9").Value = "00:00"
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AI9").Value = "Blandsumpskog"
$ws.Range("AO9").Value = "Murken klenved"
$ws.Range("AW9").Value = "Henrik Josefsson"
$ws.Range("AX9").Value = "Tommy Pettersson"
$ws.Range("AY9").Value = "Länsstyrelsen i Örebro län, inventering"

# Row 10
$ws.Range("A10").Value = 82883429
$ws.Range("B10").Value = 77506
$ws.Range("C10").Value = "Ovaliderad"
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("P10").Value = "Holmtjärnen, Vstm"
$ws.Range("Q10").Value = 489400.1580363552
$ws.Range("R10").Value = 6631499.862790138
$ws.Range("S10").Value = 25
$ws.Range("T10").Value = "Örebro"
$ws.Range("U10").Value = "Hällefors"
$ws.Range("V10").Value = "Västmanland"
$ws.Range("W10").Value = "Hjulsjö"
$ws.Range("Y10").Value = "2019-10-17"
$ws.Range("Z10").Value = "00:00"
$ws.Range("AA10").Value = "2019-11-15"
$ws.Range("AB10").Value = "00:00"
$ws.Range("AD10").Value = $false
$ws.Range("AE10").Value = $false
$ws.Range("AF10").ClearContents()
$ws.Range("AG10").Value = $false
$ws.Range("AI10").Value = "Äldre barrskog"
$ws.Range("AO10").Value = "Grankvist"
$ws.Range("AW10").Value = "Henrik Josefsson"
$ws.Range("AX10").Value = "Tommy Pettersson"
$ws.Range("AY10").Value = "Länsstyrelsen i Örebro län, inventering"

# Row 11
$ws.Range("A11").Value = 82883417
$ws.Range("B11").Value = 89356
$ws.Range("C11").Value = "Ovaliderad"
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 5447
$ws.Range("F11").Value = "Vedticka"
$ws.Range("G11").Value = "Fuscoporia viticola"
$ws.Range("H11").Value = "(Schwein.) Murrill"
$ws.Range("P11").Value = "Holmtjärnen, Vstm"
$ws.Range("Q11").Value = 489516.9242719654
$ws.Range("R11").Value = 6631461.813275921
$ws.Range("S11").Value = 25
$ws.Range("T11").Value = "Örebro"
$ws.Range("U11").Value = "Hällefors"
$ws.Range("V11").Value = "Västmanland"
$ws.Range("W11").Value = "Hjulsjö"
$ws.Range("Y11").Value = "2019-10-17"
$ws.Range("Z11").Value = "00:00"
$ws.Range("AA11").Value = "2019-11-15"
$ws.Range("AB11").Value = "00:00"
$ws.Range("AD11").Value = $false
$ws.Range("AE11").Value = $false
$ws.Range("AG11").Value = $false
$ws.Range("AI11").Value = "Barrfuktskog"
$ws.Range("AO11").Value = "Liggande granstam"
$ws.Range("AW11").Value = "Henrik Josefsson"
$ws.Range("AX11").Value = "Tommy Pettersson"
$ws.Range("AY11").Value = "Länsstyrelsen i Örebro län, inventering"
